$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 9) with a "space" keyword mapped to the existing
# "com.singleton.helix" appID, matching the pattern already used by rows
# 4/5/8 for that appID.
$ws.Range("A9").Value = "space"
$ws.Range("B9").Value = "com.singleton.helix"

# Carry over the same cell formatting (style) used by the rest of the
# data rows, e.g. row 8, onto the newly added row.
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to the newly added row, as in the source file.
$ws.Range("A9").Select()
